$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style the first cell: bold font, thin box border, centered horizontally, top vertically
$cellB1 = $ws.Range("B1")
$cellB1.Font.Bold = $true
$cellB1.HorizontalAlignment = -4108  # xlCenter
$cellB1.VerticalAlignment = -4160    # xlTop
$cellB1.Borders.LineStyle = 1        # xlContinuous

# Reuse the exact same style for A2 by copying formats only, so no extra
# (unused) style entries get created in the workbook's style table.
$cellB1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
